$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: Sending cluster ---
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "ECs"
$ws.Range("A5").Value = "ECs"
$ws.Range("A6").Value = "FAPs"
$ws.Range("A7").Value = "FAPs"
$ws.Range("A8").Value = "FAPs"
$ws.Range("A9").Value = "FAPs"
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("A11").Value = "Inflammatory-Mac"
$ws.Range("A12").Value = "Inflammatory-Mac"
$ws.Range("A13").Value = "Inflammatory-Mac"
$ws.Range("A14").Value = "MuSCs"
$ws.Range("A15").Value = "MuSCs"
$ws.Range("A16").Value = "MuSCs"
$ws.Range("A17").Value = "MuSCs"
$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("A20").Value = "Resolving-Mac"
$ws.Range("A21").Value = "Resolving-Mac"

# --- Column B: Ligand symbol ---
$ws.Range("B2").Value = "Col8a1"
$ws.Range("B3").Value = "Col8a1"
$ws.Range("B4").Value = "Col8a1"
$ws.Range("B5").Value = "Col8a1"
$ws.Range("B6").Value = "Col8a1"
$ws.Range("B7").Value = "Col8a1"
$ws.Range("B8").Value = "Col8a1"
$ws.Range("B9").Value = "Col8a1"
$ws.Range("B10").Value = "Col8a1"
$ws.Range("B11").Value = "Col8a1"
$ws.Range("B12").Value = "Col8a1"
$ws.Range("B13").Value = "Col8a1"
$ws.Range("B14").Value = "Col8a1"
$ws.Range("B15").Value = "Col8a1"
$ws.Range("B16").Value = "Col8a1"
$ws.Range("B17").Value = "Col8a1"
$ws.Range("B18").Value = "Col8a1"
$ws.Range("B19").Value = "Col8a1"
$ws.Range("B20").Value = "Col8a1"
$ws.Range("B21").Value = "Col8a1"

# --- Column C: Receptor symbol ---
$ws.Range("C2").Value = "Itga2"
$ws.Range("C3").Value = "Itga2"
$ws.Range("C4").Value = "Itga2"
$ws.Range("C5").Value = "Itga2"
$ws.Range("C6").Value = "Itga2"
$ws.Range("C7").Value = "Itga2"
$ws.Range("C8").Value = "Itga2"
$ws.Range("C9").Value = "Itga2"
$ws.Range("C10").Value = "Itga2"
$ws.Range("C11").Value = "Itga2"
$ws.Range("C12").Value = "Itga2"
$ws.Range("C13").Value = "Itga2"
$ws.Range("C14").Value = "Itga2"
$ws.Range("C15").Value = "Itga2"
$ws.Range("C16").Value = "Itga2"
$ws.Range("C17").Value = "Itga2"
$ws.Range("C18").Value = "Itga2"
$ws.Range("C19").Value = "Itga2"
$ws.Range("C20").Value = "Itga2"
$ws.Range("C21").Value = "Itga2"

# --- Column D: Target cluster ---
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("D6").Value = "ECs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("D10").Value = "ECs"
$ws.Range("D11").Value = "FAPs"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("D14").Value = "ECs"
$ws.Range("D15").Value = "FAPs"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("D18").Value = "ECs"
$ws.Range("D19").Value = "FAPs"
$ws.Range("D20").Value = "MuSCs"
$ws.Range("D21").Value = "Resolving-Mac"

# --- Column E: Ligand-expressing cells ---
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("E8").Value = 3
$ws.Range("E9").Value = 3
$ws.Range("E10").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("E14").Value = 3
$ws.Range("E15").Value = 3
$ws.Range("E16").Value = 3
$ws.Range("E17").Value = 3
$ws.Range("E18").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("E21").Value = 1

# --- Column F: Ligand detection rate ---
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("F21").Value = 0.3333333333333333

# --- Column G: Ligand average expression value ---
$ws.Range("G2").Value = 11.430265
$ws.Range("G3").Value = 11.430265
$ws.Range("G4").Value = 11.430265
$ws.Range("G5").Value = 11.430265
$ws.Range("G6").Value = 175.4151816666667
$ws.Range("G7").Value = 175.4151816666667
$ws.Range("G8").Value = 175.4151816666667
$ws.Range("G9").Value = 175.4151816666667
$ws.Range("G10").Value = 0.05623966666666667
$ws.Range("G11").Value = 0.05623966666666667
$ws.Range("G12").Value = 0.05623966666666667
$ws.Range("G13").Value = 0.05623966666666667
$ws.Range("G14").Value = 5.493340666666666
$ws.Range("G15").Value = 5.493340666666666
$ws.Range("G16").Value = 5.493340666666666
$ws.Range("G17").Value = 5.493340666666666
$ws.Range("G18").Value = 0.05570833333333333
$ws.Range("G19").Value = 0.05570833333333333
$ws.Range("G20").Value = 0.05570833333333333
$ws.Range("G21").Value = 0.05570833333333333

# --- Column H: Ligand total expression value ---
$ws.Range("H2").Value = 34.290795
$ws.Range("H3").Value = 34.290795
$ws.Range("H4").Value = 34.290795
$ws.Range("H5").Value = 34.290795
$ws.Range("H6").Value = 526.245545
$ws.Range("H7").Value = 526.245545
$ws.Range("H8").Value = 526.245545
$ws.Range("H9").Value = 526.245545
$ws.Range("H10").Value = 0.168719
$ws.Range("H11").Value = 0.168719
$ws.Range("H12").Value = 0.168719
$ws.Range("H13").Value = 0.168719
$ws.Range("H14").Value = 16.480022
$ws.Range("H15").Value = 16.480022
$ws.Range("H16").Value = 16.480022
$ws.Range("H17").Value = 16.480022
$ws.Range("H18").Value = 0.167125
$ws.Range("H19").Value = 0.167125
$ws.Range("H20").Value = 0.167125
$ws.Range("H21").Value = 0.167125

# --- Column I: Ligand derived specificity of average expression value ---
$ws.Range("I2").Value = 0.05939319992829472
$ws.Range("I3").Value = 0.05939319992829472
$ws.Range("I4").Value = 0.05939319992829472
$ws.Range("I5").Value = 0.05939319992829472
$ws.Range("I6").Value = 0.9114809634935387
$ws.Range("I7").Value = 0.9114809634935387
$ws.Range("I8").Value = 0.9114809634935387
$ws.Range("I9").Value = 0.9114809634935387
$ws.Range("I10").Value = 0.0002922288998753735
$ws.Range("I11").Value = 0.0002922288998753735
$ws.Range("I12").Value = 0.0002922288998753735
$ws.Range("I13").Value = 0.0002922288998753735
$ws.Range("I14").Value = 0.02854413965814136
$ws.Range("I15").Value = 0.02854413965814136
$ws.Range("I16").Value = 0.02854413965814136
$ws.Range("I17").Value = 0.02854413965814136
$ws.Range("I18").Value = 0.0002894680201499048
$ws.Range("I19").Value = 0.0002894680201499048
$ws.Range("I20").Value = 0.0002894680201499048
$ws.Range("I21").Value = 0.0002894680201499048

# --- Column J: Ligand derived specificity of total expression value ---
$ws.Range("J2").Value = 0.05939319992829471
$ws.Range("J3").Value = 0.05939319992829471
$ws.Range("J4").Value = 0.05939319992829471
$ws.Range("J5").Value = 0.05939319992829471
$ws.Range("J6").Value = 0.9114809634935386
$ws.Range("J7").Value = 0.9114809634935386
$ws.Range("J8").Value = 0.9114809634935386
$ws.Range("J9").Value = 0.9114809634935386
$ws.Range("J10").Value = 0.0002922288998753735
$ws.Range("J11").Value = 0.0002922288998753735
$ws.Range("J12").Value = 0.0002922288998753735
$ws.Range("J13").Value = 0.0002922288998753735
$ws.Range("J14").Value = 0.02854413965814136
$ws.Range("J15").Value = 0.02854413965814136
$ws.Range("J16").Value = 0.02854413965814136
$ws.Range("J17").Value = 0.02854413965814136
$ws.Range("J18").Value = 0.0002894680201499048
$ws.Range("J19").Value = 0.0002894680201499048
$ws.Range("J20").Value = 0.0002894680201499048
$ws.Range("J21").Value = 0.0002894680201499048

# --- Column K: Receptor-expressing cells ---
$ws.Range("K2").Value = 3
$ws.Range("K3").Value = 3
$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 3
$ws.Range("K8").Value = 3
$ws.Range("K9").Value = 1
$ws.Range("K10").Value = 3
$ws.Range("K11").Value = 3
$ws.Range("K12").Value = 3
$ws.Range("K13").Value = 1
$ws.Range("K14").Value = 3
$ws.Range("K15").Value = 3
$ws.Range("K16").Value = 3
$ws.Range("K17").Value = 1
$ws.Range("K18").Value = 3
$ws.Range("K19").Value = 3
$ws.Range("K20").Value = 3
$ws.Range("K21").Value = 1

# --- Column L: Receptor detection rate ---
$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("L10").Value = 1
$ws.Range("L11").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("L14").Value = 1
$ws.Range("L15").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("L18").Value = 1
$ws.Range("L19").Value = 1
$ws.Range("L20").Value = 1
$ws.Range("L21").Value = 0.3333333333333333

# --- Column M: Receptor average expression value ---
$ws.Range("M2").Value = 6.066605666666667
$ws.Range("M3").Value = 0.8541786666666665
$ws.Range("M4").Value = 0.1824346666666667
$ws.Range("M5").Value = 0.035773
$ws.Range("M6").Value = 6.066605666666667
$ws.Range("M7").Value = 0.8541786666666665
$ws.Range("M8").Value = 0.1824346666666667
$ws.Range("M9").Value = 0.035773
$ws.Range("M10").Value = 6.066605666666667
$ws.Range("M11").Value = 0.8541786666666665
$ws.Range("M12").Value = 0.1824346666666667
$ws.Range("M13").Value = 0.035773
$ws.Range("M14").Value = 6.066605666666667
$ws.Range("M15").Value = 0.8541786666666665
$ws.Range("M16").Value = 0.1824346666666667
$ws.Range("M17").Value = 0.035773
$ws.Range("M18").Value = 6.066605666666667
$ws.Range("M19").Value = 0.8541786666666665
$ws.Range("M20").Value = 0.1824346666666667
$ws.Range("M21").Value = 0.035773

# --- Column N: Receptor total expression value ---
$ws.Range("N2").Value = 18.199817
$ws.Range("N3").Value = 2.562536
$ws.Range("N4").Value = 0.547304
$ws.Range("N5").Value = 0.107319
$ws.Range("N6").Value = 18.199817
$ws.Range("N7").Value = 2.562536
$ws.Range("N8").Value = 0.547304
$ws.Range("N9").Value = 0.107319
$ws.Range("N10").Value = 18.199817
$ws.Range("N11").Value = 2.562536
$ws.Range("N12").Value = 0.547304
$ws.Range("N13").Value = 0.107319
$ws.Range("N14").Value = 18.199817
$ws.Range("N15").Value = 2.562536
$ws.Range("N16").Value = 0.547304
$ws.Range("N17").Value = 0.107319
$ws.Range("N18").Value = 18.199817
$ws.Range("N19").Value = 2.562536
$ws.Range("N20").Value = 0.547304
$ws.Range("N21").Value = 0.107319

# --- Column O: Receptor derived specificity of average expression value ---
$ws.Range("O2").Value = 0.849784628791665
$ws.Range("O3").Value = 0.1196497582104962
$ws.Range("O4").Value = 0.02555468148257719
$ws.Range("O5").Value = 0.005010931515261538
$ws.Range("O6").Value = 0.849784628791665
$ws.Range("O7").Value = 0.1196497582104962
$ws.Range("O8").Value = 0.02555468148257719
$ws.Range("O9").Value = 0.005010931515261538
$ws.Range("O10").Value = 0.849784628791665
$ws.Range("O11").Value = 0.1196497582104962
$ws.Range("O12").Value = 0.02555468148257719
$ws.Range("O13").Value = 0.005010931515261538
$ws.Range("O14").Value = 0.849784628791665
$ws.Range("O15").Value = 0.1196497582104962
$ws.Range("O16").Value = 0.02555468148257719
$ws.Range("O17").Value = 0.005010931515261538
$ws.Range("O18").Value = 0.849784628791665
$ws.Range("O19").Value = 0.1196497582104962
$ws.Range("O20").Value = 0.02555468148257719
$ws.Range("O21").Value = 0.005010931515261538

# --- Column P: Receptor derived specificity of total expression value ---
$ws.Range("P2").Value = 0.8497846287916652
$ws.Range("P3").Value = 0.1196497582104962
$ws.Range("P4").Value = 0.02555468148257719
$ws.Range("P5").Value = 0.005010931515261539
$ws.Range("P6").Value = 0.8497846287916652
$ws.Range("P7").Value = 0.1196497582104962
$ws.Range("P8").Value = 0.02555468148257719
$ws.Range("P9").Value = 0.005010931515261539
$ws.Range("P10").Value = 0.8497846287916652
$ws.Range("P11").Value = 0.1196497582104962
$ws.Range("P12").Value = 0.02555468148257719
$ws.Range("P13").Value = 0.005010931515261539
$ws.Range("P14").Value = 0.8497846287916652
$ws.Range("P15").Value = 0.1196497582104962
$ws.Range("P16").Value = 0.02555468148257719
$ws.Range("P17").Value = 0.005010931515261539
$ws.Range("P18").Value = 0.8497846287916652
$ws.Range("P19").Value = 0.1196497582104962
$ws.Range("P20").Value = 0.02555468148257719
$ws.Range("P21").Value = 0.005010931515261539

# --- Column Q: Edge average expression weight ---
$ws.Range("Q2").Value = 69.34291042050165
$ws.Range("Q3").Value = 9.763488517346664
$ws.Range("Q4").Value = 2.085276585186666
$ws.Range("Q5").Value = 0.408894869845
$ws.Range("Q6").Value = 1064.174735118363
$ws.Range("Q7").Value = 149.8359059891244
$ws.Range("Q8").Value = 32.00181019563111
$ws.Range("Q9").Value = 6.275127293761666
$ws.Range("Q10").Value = 0.3411838804914444
$ws.Range("Q11").Value = 0.04803872348711111
$ws.Range("Q12").Value = 0.01026006484177778
$ws.Range("Q13").Value = 0.002011861595666667
$ws.Range("Q14").Value = 33.32593161733044
$ws.Range("Q15").Value = 4.692294406199109
$ws.Range("Q16").Value = 1.002175773409778
$ws.Range("Q17").Value = 0.1965132756686666
$ws.Range("Q18").Value = 0.3379604906805556
$ws.Range("Q19").Value = 0.04758486988888888
$ws.Range("Q20").Value = 0.01016313122222222
$ws.Range("Q21").Value = 0.001992854208333333

# --- Column R: Edge total expression weight ---
$ws.Range("R2").Value = 624.0861937845149
$ws.Range("R3").Value = 87.87139665611998
$ws.Range("R4").Value = 18.76748926668
$ws.Range("R5").Value = 3.680053828604999
$ws.Range("R6").Value = 9577.572616065265
$ws.Range("R7").Value = 1348.52315390212
$ws.Range("R8").Value = 288.01629176068
$ws.Range("R9").Value = 56.476145643855
$ws.Range("R10").Value = 3.070654924423
$ws.Range("R11").Value = 0.432348511384
$ws.Range("R12").Value = 0.092340583576
$ws.Range("R13").Value = 0.018106754361
$ws.Range("R14").Value = 299.9333845559739
$ws.Range("R15").Value = 42.23064965579199
$ws.Range("R16").Value = 9.019581960687999
$ws.Range("R17").Value = 1.768619481018
$ws.Range("R18").Value = 3.041644416125
$ws.Range("R19").Value = 0.4282638289999999
$ws.Range("R20").Value = 0.091468181
$ws.Range("R21").Value = 0.017935687875

# --- Column S: Edge average expression derived specificity ---
$ws.Range("S2").Value = 0.05047142835381507
$ws.Range("S3").Value = 0.007106382010768122
$ws.Range("S4").Value = 0.001517774306398598
$ws.Range("S5").Value = 0.0002976152573129213
$ws.Range("S6").Value = 0.774562512213026
$ws.Range("S7").Value = 0.109058476895472
$ws.Range("S8").Value = 0.02329260569950994
$ws.Range("S9").Value = 0.004567368685530725
$ws.Range("S10").Value = 0.0002483316272027909
$ws.Range("S11").Value = 0.00003496511721220774
$ws.Range("S12").Value = 0.00000746781645631911
$ws.Range("S13").Value = 0.000001464339004055718
$ws.Range("S14").Value = 0.0242563711235711
$ws.Range("S15").Value = 0.003415299408423248
$ws.Range("S16").Value = 0.000729436397158002
$ws.Range("S17").Value = 0.0001430327289890072
$ws.Range("S18").Value = 0.0002459854740501451
$ws.Range("S19").Value = 0.00003463477862060714
$ws.Range("S20").Value = 0.000007397263054323053
$ws.Range("S21").Value = 0.00000145050442482952

# --- Column T: Edge total expression derived specificity ---
$ws.Range("T2").Value = 0.05047142835381507
$ws.Range("T3").Value = 0.007106382010768122
$ws.Range("T4").Value = 0.001517774306398598
$ws.Range("T5").Value = 0.0002976152573129213
$ws.Range("T6").Value = 0.774562512213026
$ws.Range("T7").Value = 0.109058476895472
$ws.Range("T8").Value = 0.02329260569950994
$ws.Range("T9").Value = 0.004567368685530725
$ws.Range("T10").Value = 0.000248331627202791
$ws.Range("T11").Value = 0.00003496511721220774
$ws.Range("T12").Value = 0.000007467816456319111
$ws.Range("T13").Value = 0.000001464339004055718
$ws.Range("T14").Value = 0.0242563711235711
$ws.Range("T15").Value = 0.003415299408423249
$ws.Range("T16").Value = 0.0007294363971580021
$ws.Range("T17").Value = 0.0001430327289890073
$ws.Range("T18").Value = 0.0002459854740501451
$ws.Range("T19").Value = 0.00003463477862060715
$ws.Range("T20").Value = 0.000007397263054323054
$ws.Range("T21").Value = 0.00000145050442482952
